# "Generate Report for Archive" - refresh the handoff status report.
#
# 1. The status text "Ready for handoff" is now stale; the localization run
#    has moved on, so every cell showing that status becomes "In Translation"
#    (Overview!E2:F3 for the zh-cn/de-de columns, and the "Status" column
#    (C2:C3) on each per-locale sheet).
# 2. The Status-related columns (Overview E:F, and column C on the zh-cn /
#    de-de sheets) are narrowed to fit the new, shorter status text.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "In Translation"

# --- Overview sheet: zh-cn (E) / de-de (F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

# Narrow the now-shorter status columns to match the refreshed content.
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("C3").Value = $statusNew
$wsZhCn.Range("C1").ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("C3").Value = $statusNew
$wsDeDe.Range("C1").ColumnWidth = 12.5
